$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disambiguate the repeated "SE / LL95 / UL95" column headers in row 1 by
# prefixing them with the risk-group they belong to (CIN2+ / CIN3+ / CANCER).

$headerUpdates = @{
    "U1"  = "CIN2+ SE immediate"
    "V1"  = "CIN2+ LL95 immediate"
    "W1"  = "CIN2+ UL95 immediate"
    "Y1"  = "CIN2+ SE 1-year"
    "Z1"  = "CIN2+ LL95 1-year"
    "AA1" = "CIN2+ UL95 1-year"
    "AC1" = "CIN2+ SE 2-year"
    "AD1" = "CIN2+ LL95 2-year"
    "AE1" = "CIN2+ UL95 2-year"
    "AG1" = "CIN2+ SE 3-year"
    "AH1" = "CIN2+ LL95 3-year"
    "AI1" = "CIN2+ UL95 3-year"
    "AK1" = "CIN2+ SE 4-year"
    "AL1" = "CIN2+ LL95 4-year"
    "AM1" = "CIN2+ UL95 4-year"
    "AO1" = "CIN2+ SE 5-year"
    "AP1" = "CIN2+ LL95 5-year"
    "AQ1" = "CIN2+ UL95 5-year"
    "AS1" = "CIN3+ SE immediate"
    "AT1" = "CIN3+ LL95 immediate"
    "AU1" = "CIN3+ UL95 immediate"
    "AW1" = "CIN3+ SE 1-year"
    "AX1" = "CIN3+ LL95 1-year"
    "AY1" = "CIN3+ UL95 1-year"
    "BA1" = "CIN3+ SE 2-year"
    "BB1" = "CIN3+ LL95 2-year"
    "BC1" = "CIN3+ UL95 2-year"
    "BE1" = "CIN3+ SE 3-year"
    "BF1" = "CIN3+ LL95 3-year"
    "BG1" = "CIN3+ UL95 3-year"
    "BI1" = "CIN3+ SE 4-year"
    "BJ1" = "CIN3+ LL95 4-year"
    "BK1" = "CIN3+ UL95 4-year"
    "BM1" = "CIN3+ SE 5-year"
    "BN1" = "CIN3+ LL95 5-year"
    "BO1" = "CIN3+ UL95 5-year"
    "BQ1" = "CANCER SE immediate"
    "BR1" = "CANCER LL95 immediate"
    "BS1" = "CANCER UL95 immediate"
    "BU1" = "CANCER E 1-year"
    "BV1" = "CANCER LL95 1-year"
    "BW1" = "CANCER UL95 1-year"
    "BY1" = "CANCER SE 2-year"
    "BZ1" = "CANCER LL95 2-year"
    "CA1" = "CANCER UL95 2-year"
    "CC1" = "CANCER SE 3-year"
    "CD1" = "CANCER LL95 3-year"
    "CE1" = "CANCER UL95 3-year"
    "CG1" = "CANCER SE 4-year"
    "CH1" = "CANCER LL95 4-year"
    "CI1" = "CANCER UL95 4-year"
    "CK1" = "CANCER SE 5-year"
    "CL1" = "CANCER LL95 5-year"
    "CM1" = "CANCER UL95 5-year"
}

foreach ($addr in $headerUpdates.Keys) {
    $ws.Range($addr).Value = $headerUpdates[$addr]
}

# Move the frozen-pane view / selection to match the updated layout.
$ws.Range("F1:CP1").Select()
$excel.ActiveWindow.ScrollColumn = 11
